# Preparing for a staging release.docx - apply commit edits:
#  1) Add "CloudSetupApiSyncSampleSupport" to the assembly/file-version paragraph.
#  2) Add a new paragraph about changing the CloudSetupSdkSyncSample version.
#  3) Add a new "Exit Visual Studio." paragraph after the testing step.
#  4) Move the _GoBack bookmark to the start of the "Commit and push..." paragraph.

$d = $word.ActiveDocument

# --- 1) "Change the CloudApiPublic Assembly and File Versions..." paragraph ---
$pAssembly = $d.Paragraphs(4)
$pAssembly.Range.Find.Execute(
    " Assembly and File Versions to the current release.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and CloudSetupApiSyncSampleSupport Assembly and File Versions to the current release.  ",
    2
)

# --- 2) Insert a new paragraph right after it, describing the CloudSetupSdkSyncSample version change ---
$pAssembly = $d.Paragraphs(4)
$rEnd = $pAssembly.Range.Duplicate
$rEnd.Collapse(0)
$rEnd.InsertParagraphAfter()
$pNewVersion = $d.Paragraphs(5)
$pNewVersion.Range.Text = "Change the CloudSetupSdkSyncSample version.  Click Project Assistant and change the application version.  E.g., 0.01.0002."

# --- 3) Insert "Exit Visual Studio." after "Test the application..." paragraph ---
$pTest = $d.Paragraphs(6)
$rTestEnd = $pTest.Range.Duplicate
$rTestEnd.Collapse(0)
$rTestEnd.InsertParagraphAfter()
$pExit = $d.Paragraphs(7)
$pExit.Range.Text = "Exit Visual Studio."

# --- 4) Move the hidden _GoBack bookmark to the start of the "Commit and push..." paragraph ---
$pCommit = $d.Paragraphs(8)
$rBookmark = $pCommit.Range.Duplicate
$rBookmark.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rBookmark)
